$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated values (Coin name, Link, Price, Volume(1h)) per the source diff.
# Price (column D) cells get an explicit Text format first so that numeric-
# looking strings (e.g. "232.23") are not silently reinterpreted by Excel as
# numbers, matching the original inline-string cell content.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '41.907.22'
$ws.Range("E2").Value = '  +1.13%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.235.77'
$ws.Range("E3").Value = '  -0.37%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.23'
$ws.Range("E5").Value = '  +0.99%  '
$ws.Range("E6").Value = '  -1.91%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '60.51'
$ws.Range("E7").Value = '  -6.65%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.405'
$ws.Range("E9").Value = '  -0.54%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '57.99'
$ws.Range("E10").Value = '  -3.89%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0906'
$ws.Range("E11").Value = '  -0.03%  '
$ws.Range("E12").Value = '  -0.49%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.570.25'
$ws.Range("E13").Value = '  -0.27%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.56'
$ws.Range("E14").Value = '  -3.51%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '22.83'
$ws.Range("E15").Value = '  +2.41%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.66'
$ws.Range("E16").Value = '  +0.59%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.805'
$ws.Range("E17").Value = '  -2.72%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.246.72'
$ws.Range("E18").Value = '  +0.24%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '41.897.17'
$ws.Range("E19").Value = '  +1.37%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0916'
$ws.Range("E20").Value = '  -0.14%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.65'
$ws.Range("E21").Value = '  -1.74%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.16'
$ws.Range("E22").Value = '  +0.25%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '248.66'
$ws.Range("E23").Value = '  -1.43%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.00'
$ws.Range("E24").Value = '  -0.09%  '
$ws.Range("B25").Value = 'PancakeSwap'
$ws.Range("C25").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.38'
$ws.Range("E25").Value = '  -0.88%  '
$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.40'
$ws.Range("E26").Value = '  +3.31%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.77'
$ws.Range("E27").Value = '  +0.16%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '169.31'
$ws.Range("E28").Value = '  -2.00%  '
$ws.Range("E29").Value = '  -2.12%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.95'
$ws.Range("E31").Value = '  -1.14%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.56'
$ws.Range("E32").Value = '  -9.82%  '
$ws.Range("E33").Value = '  -1.54%  '
$ws.Range("E34").Value = '  +3.03%  '
$ws.Range("E35").Value = '  -0.41%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0658'
$ws.Range("E36").Value = '  +4.07%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.62'
$ws.Range("E37").Value = '  -8.09%  '
$ws.Range("E38").Value = '  -1.05%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.64'
$ws.Range("E39").Value = '  -5.42%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.000248'
$ws.Range("E40").Value = '  +6.22%  '
$ws.Range("E41").Value = '  -0.12%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0240'
$ws.Range("E42").Value = '  +1.51%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.68'
$ws.Range("E43").Value = '  -1.13%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.22'
$ws.Range("E44").Value = '  -0.10%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '99.39'
$ws.Range("E45").Value = '  -2.25%  '
$ws.Range("B46").Value = 'FTXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.48'
$ws.Range("E46").Value = '  -8.06%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0965'
$ws.Range("E47").Value = '  +2.55%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.478.21'
$ws.Range("E48").Value = '  -2.29%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '16.74'
$ws.Range("E49").Value = '  -6.51%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.30'
$ws.Range("E50").Value = '  +8.91%  '
$ws.Range("B51").Value = 'ARBITRUM'
$ws.Range("C51").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.09'
